$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 551
$ws.Range("I2").Value = 1510
$ws.Range("J2").Value = 6376
$ws.Range("K2").Value = 33
$ws.Range("L2").Value = 1707
$ws.Range("M2").Value = 96
$ws.Range("N2").Value = 1163
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 22
$ws.Range("R2").Value = 90
$ws.Range("S2").Value = 697
$ws.Range("T2").Value = 1128
$ws.Range("U2").Value = 91
$ws.Range("V2").Value = 9739
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 9657
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 157
$ws.Range("AA2").Value = 50
